# Update minimal example (besprechung shb jf)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A (Var.Name) updates ---
$ws.Range("A2").Value = "ID"
$ws.Range("A6").Value = "skala1_item1"
$ws.Range("A7").Value = "skala1_item2"
$ws.Range("A8").Value = "skala1_item3"
$ws.Range("A11").Value = "pv_1"
$ws.Range("A12").Value = "pv_2"
$ws.Range("A13").Value = "pv_3"
$ws.Range("A14").Value = "pv_4"
$ws.Range("A15").Value = "pv_5"
$ws.Range("A16").Value = "pvkat_1"
$ws.Range("A17").Value = "pvkat_2"
$ws.Range("A18").Value = "pvkat_3"
$ws.Range("A19").Value = "pvkat_4"
$ws.Range("A20").Value = "pvkat_5"

# --- Column D (LabelSH) and Column H (Titel) updates ---
$ws.Range("D9").Value = "Skala: Likert-Skalenwert"
$ws.Range("H9").Value = "Skala: Likert-Skalenwert"

$ws.Range("D10").Value = "IMPUTATION 1: plausible value"
$ws.Range("H10").Value = "IMPUTATION 1: plausible value"
$ws.Range("D11").Value = "IMPUTATION 1: plausible value"
$ws.Range("H11").Value = "IMPUTATION 1: plausible value"

$ws.Range("D12").Value = "IMPUTATION 2: plausible value"
$ws.Range("H12").Value = "IMPUTATION 2: plausible value"

$ws.Range("D13").Value = "IMPUTATION 3: plausible value"
$ws.Range("H13").Value = "IMPUTATION 3: plausible value"

$ws.Range("D14").Value = "IMPUTATION 4: plausible value"
$ws.Range("H14").Value = "IMPUTATION 4: plausible value"

$ws.Range("D15").Value = "IMPUTATION 5: plausible value"
$ws.Range("H15").Value = "IMPUTATION 5: plausible value"

$ws.Range("D16").Value = "IMPUTATION 1: Kompetenzstufe des plausible value"
$ws.Range("H16").Value = "IMPUTATION 1: Kompetenzstufe des plausible value"

$ws.Range("D17").Value = "IMPUTATION 2: Kompetenzstufe des plausible value"
$ws.Range("H17").Value = "IMPUTATION 2: Kompetenzstufe des plausible value"

$ws.Range("D18").Value = "IMPUTATION 3: Kompetenzstufe des plausible value"
$ws.Range("H18").Value = "IMPUTATION 3: Kompetenzstufe des plausible value"

$ws.Range("D19").Value = "IMPUTATION 4: Kompetenzstufe des plausible value"
$ws.Range("H19").Value = "IMPUTATION 4: Kompetenzstufe des plausible value"

$ws.Range("D20").Value = "IMPUTATION 4: Kompetenzstufe des plausible value"
$ws.Range("H20").Value = "IMPUTATION 4: Kompetenzstufe des plausible value"

$ws.Range("D24").Value = "Skala: fake-skala"
$ws.Range("H24").Value = "Skala: fake-skala"
